# Ajuste na pontuação das tabelas
# Update the COP_B_PTS (column I) values on the active sheet rows 2-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 133.33
    3  = 116.67
    4  = 100
    5  = 55.56
    6  = 53.33
    7  = 38.89
    8  = 16.67
    9  = 14.44
    10 = 11.11
    11 = 11.11
    12 = 8.890000000000001
    13 = 8.890000000000001
    14 = 5.56
    15 = 5.56
    16 = 3.33
    17 = 3.33
    18 = 3.33
    19 = 3.33
    20 = 3.33
    21 = 3.33
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 9).Value = $newValues[$row]
}
